# 02_Middleware.pptx — "updating for core 2"
#
# The underlying OOXML diff is dominated by PowerPoint's own file/rId
# renumbering that happens as a side effect of removing two slides
# ("Architecture" and "Processing HTTP Messages") from the deck. The real,
# semantic edit is:
#   1. Delete slide 2 ("Architecture").
#   2. Delete the (now) slide 2 ("Processing HTTP Messages").
#   3. Refresh the cached handout-master date field to the new save date.
#
# (The cached slide-number field baked into the "Tips" slide's notes page
#  would also read "13" instead of "15" after the two deletions above, but
#  that value is recomputed by PowerPoint automatically whenever the notes
#  page/field is redrawn - it is not something a user edits by hand.)

$p = $ppt.ActivePresentation

# --- 1 & 2. Remove the two slides that were dropped from the deck -------
# "Architecture" sits at index 2, "Processing HTTP Messages" right after it
# at index 3. Deleting index 2 twice removes both, because the second slide
# slides up into index 2 once the first is gone.
$p.Slides.Item(2).Delete()
$p.Slides.Item(2).Delete()

# --- 3. Update the cached date field on the handout master ---------------
$hm = $ppt.ActivePresentation.HandoutMaster
$hf = $hm.HeadersFooters
$hf.DateAndTime.Text = "3/8/2018"
